$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is the "pdf_online" record: update its name/description
# from "pdf online" to "pdf online test"
$ws.Range("C4").Value = "pdf online test"
$ws.Range("D4").Value = "pdf online test"

# Update the active selection on the sheet
$ws.Range("C10").Select()
